# Ütemterv (schedule) sheet: extend the Gantt grid down to row 22, draw the
# three vertical "milestone" divider lines in columns C, G and L, and label
# them M1 / M2 / M3 on the new bottom row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeRight = 10
$xlContinuous = 1
$xlRight = -4152

# Vertical milestone divider lines: a thin border on the right edge of every
# cell in columns C, G and L from row 3 down to row 21.
foreach ($col in @("C", "G", "L")) {
    $line = $ws.Range($col + "3:" + $col + "21")
    $line.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
}

# Milestone labels on the new row 22, bold + right aligned, with the same
# right border as the rest of the divider line above them.
$milestoneAddrs = @("C22", "G22", "L22")
$milestoneText = @("M1", "M2", "M3")
for ($i = 0; $i -lt $milestoneAddrs.Length; $i++) {
    $cell = $ws.Range($milestoneAddrs[$i])
    $cell.Value = $milestoneText[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlRight
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
}

# Match the author's final selection.
$null = $ws.Range("H28").Select()
